$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2: image_0018.jpg -> image_0017.jpg
$ws.Range("B2").Value = "image_0017.jpg"

# Update C2: dog -> invalid_stamp
$ws.Range("C2").Value = "invalid_stamp"

# Update G2: invalid symbol or invalid stamp -> invalid stamp
$ws.Range("G2").Value = "invalid stamp"

# Update H2: numeric 237.3010746068453 -> inline string "nan"
$ws.Range("H2").Value = "nan"
